$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("PFS caba").Name = "PFS_caba"
$wb.Worksheets.Item("PFS mito").Name = "PFS_mito"
$wb.Worksheets.Item("OS caba").Name = "OS_caba"
$wb.Worksheets.Item("OS mito").Name = "OS_mito"
